# DOMA-1155 regenerate excel templates
# Insert a "Unit" column after "Address" and split the single "Reading"
# column into four per-tariff reading columns, pushing "Contact"/"Source"
# to the right. Net effect: header/data grid grows from A1:H3 to A1:L3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new column at C ("Unit") - shifts old C..H (Service..Source) to D..I
$ws.Columns.Item(3).Insert()

# Insert three more new columns at G (after the new D..F: Service, Meter number, Place)
# so the single "Reading" column (currently at G) becomes four tariff columns G:J,
# and old Contact/Source (H:I) shift to K:L
$ws.Range("G1:I1").EntireColumn.Insert()

# ---- Row 1: headers ----
$ws.Range("A1").Value = "Reading date"
$ws.Range("B1").Value = "Address"
$ws.Range("C1").Value = "Unit"
$ws.Range("D1").Value = "Service"
$ws.Range("E1").Value = "Meter number"
$ws.Range("F1").Value = "Place"
$ws.Range("G1").Value = "Reading from tariff №1"
$ws.Range("H1").Value = "Reading from tariff №2"
$ws.Range("I1").Value = "Reading from tariff №3"
$ws.Range("J1").Value = "Reading from tariff №4"
$ws.Range("K1").Value = "Contact"
$ws.Range("L1").Value = "Source"

# ---- Row 2: template placeholders for item i ----
$ws.Range("A2").Value = "{d.meter[i].date}"
$ws.Range("B2").Value = "{d.meter[i].address}"
$ws.Range("C2").Value = "{d.meter[i].unitName}"
$ws.Range("D2").Value = "{d.meter[i].resource}"
$ws.Range("E2").Value = "{d.meter[i].number}"
$ws.Range("F2").Value = "{d.meter[i].place}"
$ws.Range("G2").Value = "{d.meter[i].value1}"
$ws.Range("H2").Value = "{d.meter[i].value2}"
$ws.Range("I2").Value = "{d.meter[i].value3}"
$ws.Range("J2").Value = "{d.meter[i].value4}"
$ws.Range("K2").Value = "{d.meter[i].clientName}"
$ws.Range("L2").Value = "{d.meter[i].source}"

# ---- Row 3: template placeholders for item i + 1 ----
$ws.Range("A3").Value = "{d.meter[i + 1].date}"
$ws.Range("B3").Value = "{d.meter[i + 1].address}"
$ws.Range("C3").Value = "{d.meter[i + 1].unitName}"
$ws.Range("D3").Value = "{d.meter[i + 1].resource}"
$ws.Range("E3").Value = "{d.meter[i + 1].number}"
$ws.Range("F3").Value = "{d.meter[i + 1].place}"
$ws.Range("G3").Value = "{d.meter[i + 1].value1}"
$ws.Range("H3").Value = "{d.meter[i + 1].value2}"
$ws.Range("I3").Value = "{d.meter[i + 1].value3}"
$ws.Range("J3").Value = "{d.meter[i + 1].value4}"
$ws.Range("K3").Value = "{d.meter[i + 1].clientName}"
$ws.Range("L3").Value = "{d.meter[i + 1].source}"
